$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style) from the last existing header cell (G1) to new header cell (H1)
$ws.Cells.Item(1, 7).Copy()
$ws.Cells.Item(1, 8).PasteSpecial(-4122)
$ws.Cells.Item(1, 8).Value = "Save"

# Fill in the new "Save" column values for rows 2-5
$ws.Cells.Item(2, 8).Value = 1
$ws.Cells.Item(3, 8).Value = 0
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(5, 8).Value = 1
